{"js": "// The worksheet table has 20 rows x 5 columns of simple arithmetic\n// expressions (e.g. \"91-5=\"). This edit replaces each cell's expression\n// with a newly generated one, in row-major (reading) order, matching the\n// document's existing cell order exactly.\nconst newValues = [\n  [\"92-53=\", \"35-23=\", \"0+53=\", \"75-60=\", \"75-68=\"],\n  [\"21+53=\", \"91-45=\", \"11+70=\", \"82-69=\", \"67-0=\"],\n  [\"54+15=\", \"46-5=\", \"51-49=\", \"60-52=\", \"68+18=\"],\n  [\"90-40=\", \"69-65=\", \"25-0=\", \"93-57=\", \"62-21=\"],\n  [\"54-11=\", \"57-24=\", \"53+24=\", \"31+38=\", \"5+62=\"],\n  [\"68-36=\", \"75-7=\", \"46+5=\", \"8+16=\", \"69-55=\"],\n  [\"49-15=\", \"82-79=\", \"54+33=\", \"8+8=\", \"11+29=\"],\n  [\"97-48=\", \"37-15=\", \"70-39=\", \"9+29=\", \"30+23=\"],\n  [\"42+47=\", \"94-88=\", \"65-65=\", \"16+63=\", \"1+76=\"],\n  [\"59+0=\", \"32+57=\", \"58-44=\", \"49+19=\", \"53+43=\"],\n  [\"1+9=\", \"14+14=\", \"24+62=\", \"78-6=\", \"87-16=\"],\n  [\"51-1=\", \"82-66=\", \"84-45=\", \"62-44=\", \"30+57=\"],\n  [\"81-48=\", \"72-31=\", \"78-11=\", \"95-52=\", \"52+1=\"],\n  [\"40-15=\", \"28+8=\", \"80-76=\", \"1+17=\", \"11-7=\"],\n  [\"94-69=\", \"3+91=\", \"60+22=\", \"67-7=\", \"33+0=\"],\n  [\"4+90=\", \"42+13=\", \"26+66=\", \"54+45=\", \"1+90=\"],\n  [\"13+10=\", \"57-56=\", \"24+10=\", \"15+16=\", \"26+5=\"],\n  [\"65-63=\", \"36+49=\", \"57+22=\", \"29+37=\", \"10+65=\"],\n  [\"15-11=\", \"66+27=\", \"41-28=\", \"67-37=\", \"0+7=\"],\n  [\"73-55=\", \"55-37=\", \"49+3=\", \"50+38=\", \"36-15=\"]\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (let r = 0; r < rows.items.length; r++) {\n  const cells = rows.items[r].cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (let c = 0; c < cells.items.length; c++) {\n    const newText = newValues[r][c];\n    if (newText === undefined) continue;\n    // Writing through the cell body's paragraph keeps the existing run\n    // formatting (font/size) intact - only the text content changes.\n    const cell = cells.items[c];\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n\n    if (paragraphs.items.length > 0) {\n      paragraphs.items[0].insertText(newText, Word.InsertLocation.replace);\n    } else {\n      cell.body.insertText(newText, Word.InsertLocation.replace);\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document's single table holds a 20x5 grid of simple arithmetic\n# expressions (e.g. \"91-5=\"). This script replaces every cell's expression\n# with a newly generated one, walking the table in row-major (reading)\n# order so each cell lines up with the matching entry in $newValues.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"92-53=\", \"35-23=\", \"0+53=\", \"75-60=\", \"75-68=\"),\n    @(\"21+53=\", \"91-45=\", \"11+70=\", \"82-69=\", \"67-0=\"),\n    @(\"54+15=\", \"46-5=\", \"51-49=\", \"60-52=\", \"68+18=\"),\n    @(\"90-40=\", \"69-65=\", \"25-0=\", \"93-57=\", \"62-21=\"),\n    @(\"54-11=\", \"57-24=\", \"53+24=\", \"31+38=\", \"5+62=\"),\n    @(\"68-36=\", \"75-7=\", \"46+5=\", \"8+16=\", \"69-55=\"),\n    @(\"49-15=\", \"82-79=\", \"54+33=\", \"8+8=\", \"11+29=\"),\n    @(\"97-48=\", \"37-15=\", \"70-39=\", \"9+29=\", \"30+23=\"),\n    @(\"42+47=\", \"94-88=\", \"65-65=\", \"16+63=\", \"1+76=\"),\n    @(\"59+0=\", \"32+57=\", \"58-44=\", \"49+19=\", \"53+43=\"),\n    @(\"1+9=\", \"14+14=\", \"24+62=\", \"78-6=\", \"87-16=\"),\n    @(\"51-1=\", \"82-66=\", \"84-45=\", \"62-44=\", \"30+57=\"),\n    @(\"81-48=\", \"72-31=\", \"78-11=\", \"95-52=\", \"52+1=\"),\n    @(\"40-15=\", \"28+8=\", \"80-76=\", \"1+17=\", \"11-7=\"),\n    @(\"94-69=\", \"3+91=\", \"60+22=\", \"67-7=\", \"33+0=\"),\n    @(\"4+90=\", \"42+13=\", \"26+66=\", \"54+45=\", \"1+90=\"),\n    @(\"13+10=\", \"57-56=\", \"24+10=\", \"15+16=\", \"26+5=\"),\n    @(\"65-63=\", \"36+49=\", \"57+22=\", \"29+37=\", \"10+65=\"),\n    @(\"15-11=\", \"66+27=\", \"41-28=\", \"67-37=\", \"0+7=\"),\n    @(\"73-55=\", \"55-37=\", \"49+3=\", \"50+38=\", \"36-15=\")\n)\n\n$rowCount = [Math]::Min($t.Rows.Count, $newValues.Count)\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    $rowValues = $newValues[$r - 1]\n    $colCount = [Math]::Min($t.Columns.Count, $rowValues.Count)\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        # Setting Range.Text replaces only the cell's text content; the\n        # existing run formatting (font/size) on the cell's paragraph mark\n        # is preserved by Word.\n        $cell.Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
